$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 67.77778000000001
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws.Range("H92").Value = 1362.2778
$ws.Range("I92").Value = 1514.8
$ws.Range("K92").Value = 1514.8
$ws.Range("M92").Value = -266.8

$ws.Range("H96").Value = 1749.381
$ws.Range("I96").Value = 672.375
$ws.Range("K96").Value = 2017.125
$ws.Range("M96").Value = -644.125

$ws.Range("H104").Value = 211.66667
$ws.Range("I104").Value = 211.66667
$ws.Range("K104").Value = 635.00001
$ws.Range("M104").Value = 1111.99999

$ws.Range("H132").Value = 9202.25
$ws.Range("J132").Value = 598
$ws.Range("L132").Value = 1794
$ws.Range("N132").Value = -6854

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1042580.3
$ws.Range("I74").Value = 1247385.1
$ws.Range("J74").Value = 18556
$ws.Range("K74").Value = 1247385.1
$ws.Range("L74").Value = 18556
$ws.Range("M74").Value = -1246511.1
$ws.Range("N74").Value = -20304

$ws.Range("H77").Value = 1042580.3
$ws.Range("I77").Value = 1247385.1
$ws.Range("J77").Value = 18556
$ws.Range("K77").Value = 6236925.5
$ws.Range("L77").Value = 92780
$ws.Range("M77").Value = -6232557.5
$ws.Range("N77").Value = -101516

$ws.Range("H122").Value = 2991.9
$ws.Range("I122").Value = 2612.75
$ws.Range("J122").Value = 3244.6667
$ws.Range("K122").Value = 7838.25
$ws.Range("L122").Value = 9734.000100000001
$ws.Range("M122").Value = -5388.25
$ws.Range("N122").Value = -14634.0001

$ws.Range("H132").Value = 3576.5273
$ws.Range("I132").Value = 1991.8889
$ws.Range("J132").Value = 6579
$ws.Range("K132").Value = 5975.6667
$ws.Range("L132").Value = 19737
$ws.Range("M132").Value = -3445.6667
$ws.Range("N132").Value = -24797

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 54871.875
$ws.Range("I20").Value = 59111.406
$ws.Range("K20").Value = 59111.406
$ws.Range("M20").Value = -58864.406

$ws.Range("H22").Value = 621.4286
$ws.Range("I22").Value = 366.66666
$ws.Range("J22").Value = 812.5
$ws.Range("K22").Value = 366.66666
$ws.Range("L22").Value = 812.5
$ws.Range("M22").Value = -193.66666
$ws.Range("N22").Value = -1158.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1747.6666
$ws.Range("I22").Value = 1087.5454
$ws.Range("K22").Value = 1087.5454
$ws.Range("M22").Value = -737.5454

$ws.Range("H31").Value = 8751679
$ws.Range("I31").Value = 3078373.8
$ws.Range("K31").Value = 3078373.8
$ws.Range("M31").Value = -3078078.8

$ws.Range("H34").Value = 8751679
$ws.Range("I34").Value = 3078373.8
$ws.Range("K34").Value = 3078373.8
$ws.Range("M34").Value = -3078171.8

$ws.Range("H58").Value = 6097000
$ws.Range("I58").Value = 11906400
$ws.Range("J58").Value = 2560844
$ws.Range("K58").Value = 11906400
$ws.Range("L58").Value = 2560844
$ws.Range("M58").Value = -11906197
$ws.Range("N58").Value = -2561250

$ws.Range("H86").Value = 43136.08
$ws.Range("I86").Value = 136403.67
$ws.Range("J86").Value = 13683.158
$ws.Range("K86").Value = 136403.67
$ws.Range("L86").Value = 13683.158
$ws.Range("M86").Value = -135280.67
$ws.Range("N86").Value = -15929.158

$ws.Range("H89").Value = 43136.08
$ws.Range("I89").Value = 136403.67
$ws.Range("J89").Value = 13683.158
$ws.Range("K89").Value = 682018.3500000001
$ws.Range("L89").Value = 68415.78999999999
$ws.Range("M89").Value = -676402.3500000001
$ws.Range("N89").Value = -79647.78999999999

$ws.Range("H132").Value = 5598.6
$ws.Range("I132").Value = 6079.4287
$ws.Range("K132").Value = 18238.2861
$ws.Range("M132").Value = -15708.2861

$ws.Range("H134").Value = 3009.3333
$ws.Range("I134").Value = 3126.158
$ws.Range("J134").Value = 1899.5
$ws.Range("K134").Value = 9378.474
$ws.Range("L134").Value = 5698.5
$ws.Range("M134").Value = -6843.474
$ws.Range("N134").Value = -10768.5

$ws.Range("H136").Value = 6097000
$ws.Range("I136").Value = 11906400
$ws.Range("J136").Value = 2560844
$ws.Range("K136").Value = 35719200
$ws.Range("L136").Value = 7682532
$ws.Range("M136").Value = -35716650
$ws.Range("N136").Value = -7687632

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 45714.285
$ws.Range("J37").Value = 45714.285
$ws.Range("L37").Value = 137142.855
$ws.Range("N37").Value = -137366.855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 999.7273
$ws.Range("I102").Value = 899.7
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 899.7
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = 722.3
$ws.Range("N102").Value = -5244

$ws.Range("H103").Value = 37151
$ws.Range("J103").Value = 37151
$ws.Range("L103").Value = 37151
$ws.Range("N103").Value = -39495

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7049.154
$ws.Range("I40").Value = 6767.625
$ws.Range("K40").Value = 6767.625
$ws.Range("M40").Value = -6631.625

$ws.Range("H61").Value = 12406.818
$ws.Range("I61").Value = 9747.5
$ws.Range("K61").Value = 9747.5
$ws.Range("M61").Value = -9545.5

$ws.Range("H93").Value = 3942.6365
$ws.Range("J93").Value = 7367.25
$ws.Range("L93").Value = 7367.25
$ws.Range("N93").Value = -9863.25

$ws.Range("H113").Value = 12406.818
$ws.Range("I113").Value = 9747.5
$ws.Range("K113").Value = 9747.5
$ws.Range("M113").Value = -7577.5

$ws.Range("H132").Value = 2384497
$ws.Range("I132").Value = 3791466
$ws.Range("J132").Value = 3472.6155
$ws.Range("K132").Value = 11374398
$ws.Range("L132").Value = 10417.8465
$ws.Range("M132").Value = -11371868
$ws.Range("N132").Value = -15477.8465

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 43042.406
$ws.Range("I122").Value = 1319.1578
$ws.Range("J122").Value = 142135.12
$ws.Range("K122").Value = 3957.4734
$ws.Range("L122").Value = 426405.36
$ws.Range("M122").Value = -1507.4734
$ws.Range("N122").Value = -431305.36

$ws.Range("H132").Value = 3789430.2
$ws.Range("I132").Value = 4387348
$ws.Range("J132").Value = 2620
$ws.Range("K132").Value = 13162044
$ws.Range("L132").Value = 7860
$ws.Range("M132").Value = -13159514
$ws.Range("N132").Value = -12920
